$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.460039973258972
$ws.Range("B1").Value = 2.030557870864868
$ws.Range("C1").Value = 2.517162799835205
$ws.Range("D1").Value = 4.7967848777771
$ws.Range("E1").Value = 0.7974292635917664
